$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.839.47"
$ws.Range("E2").Value = "  -4.24%  "
$ws.Range("D3").Value = "3.682.41"
$ws.Range("E3").Value = "  -1.72%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "583.88"
$ws.Range("E5").Value = "  -4.79%  "
$ws.Range("D6").Value = "167.95"
$ws.Range("E6").Value = "  -6.01%  "
$ws.Range("D7").Value = "3.679.38"
$ws.Range("E7").Value = "  -1.84%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").Value = "0.511"
$ws.Range("E9").Value = "  -3.60%  "
$ws.Range("D10").Value = "0.154"
$ws.Range("E10").Value = "  -7.52%  "
$ws.Range("D11").Value = "6.05"
$ws.Range("E11").Value = "  -8.17%  "
$ws.Range("D12").Value = "0.449"
$ws.Range("E12").Value = "  -6.89%  "
$ws.Range("D13").Value = "36.79"
$ws.Range("E13").Value = "  -8.33%  "
$ws.Range("D14").Value = "0.0000237"
$ws.Range("E14").Value = "  -6.89%  "
$ws.Range("D15").Value = "4.318.63"
$ws.Range("E15").Value = "  -1.09%  "
$ws.Range("D16").Value = "3.680.54"
$ws.Range("E16").Value = "  -1.70%  "
$ws.Range("D17").Value = "66.872.37"
$ws.Range("E17").Value = "  -4.21%  "
$ws.Range("D18").Value = "0.114"
$ws.Range("E18").Value = "  -5.36%  "
$ws.Range("D19").Value = "6.97"
$ws.Range("E19").Value = "  -6.61%  "
$ws.Range("D20").Value = "15.69"
$ws.Range("E20").Value = "  -4.47%  "
$ws.Range("D21").Value = "476.17"
$ws.Range("E21").Value = "  -5.28%  "
$ws.Range("D22").Value = "8.80"
$ws.Range("E22").Value = "  -4.14%  "
$ws.Range("D23").Value = "0.700"
$ws.Range("E23").Value = "  -3.08%  "
$ws.Range("D24").Value = "82.40"
$ws.Range("E24").Value = "  -4.29%  "
$ws.Range("D25").Value = "2.28"
$ws.Range("E25").Value = "  -12.96%  "
$ws.Range("D26").Value = "0.0000133"
$ws.Range("E26").Value = "  -1.73%  "
$ws.Range("D27").Value = "11.91"
$ws.Range("E27").Value = "  -8.03%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "9.93"
$ws.Range("E29").Value = "  -12.68%  "
$ws.Range("D30").Value = "2.84"
$ws.Range("E30").Value = "  -2.57%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "2.33"
$ws.Range("E31").Value = "  -5.77%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "31.49"
$ws.Range("E32").Value = "  +3.30%  "
$ws.Range("D33").Value = "7.51"
$ws.Range("E33").Value = "  -7.70%  "
$ws.Range("D34").Value = "0.105"
$ws.Range("E34").Value = "  -7.07%  "
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("D36").Value = "0.980"
$ws.Range("E36").Value = "  -6.67%  "
$ws.Range("D37").Value = "0.131"
$ws.Range("E37").Value = "  -4.40%  "
$ws.Range("D38").Value = "5.55"
$ws.Range("E38").Value = "  -9.51%  "
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").Value = "0.315"
$ws.Range("E39").Value = "  -10.35%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "435.40"
$ws.Range("E40").Value = "  -2.42%  "
$ws.Range("D41").Value = "48.44"
$ws.Range("E41").Value = "  -2.61%  "
$ws.Range("D42").Value = "1.94"
$ws.Range("E42").Value = "  -6.03%  "
$ws.Range("D43").Value = "2.74"
$ws.Range("E43").Value = "  -11.25%  "
$ws.Range("D44").Value = "8.08"
$ws.Range("E44").Value = "  -5.68%  "
$ws.Range("D45").Value = "40.38"
$ws.Range("E45").Value = "  -10.87%  "
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Value = "140.17"
$ws.Range("E46").Value = "  +1.25%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.748.39"
$ws.Range("E48").Value = "  -6.88%  "
$ws.Range("D49").Value = "0.0339"
$ws.Range("E49").Value = "  -6.04%  "
$ws.Range("D50").Value = "25.16"
$ws.Range("E50").Value = "  -7.34%  "
$ws.Range("D51").Value = "22.44"
$ws.Range("E51").Value = "  +4.83%  "
